$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column N (year 2022) data, mirroring the formatting of column M ---

# N2: blank cell with the same bottom-border-only formatting as M2
$ws.Range("M2").Copy() | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null

# N3: header year value (2022), formatted like M3 (bold, right aligned, bordered)
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null
$ws.Range("N3").Value = 2022

# N4: data value, formatted like M4
$ws.Range("M4").Copy() | Out-Null
$ws.Range("N4").PasteSpecial(-4122) | Out-Null
$ws.Range("N4").Value = 6333

# N5: data value, formatted like M5
$ws.Range("M5").Copy() | Out-Null
$ws.Range("N5").PasteSpecial(-4122) | Out-Null
$ws.Range("N5").Value = 82675

# N6: data value, formatted like M6 (bottom border, thick)
$ws.Range("M6").Copy() | Out-Null
$ws.Range("N6").PasteSpecial(-4122) | Out-Null
$ws.Range("N6").Value = 300853

# Clear clipboard marching ants / selection artifacts, then set the active cell to N2
$ws.Range("N2").Select() | Out-Null
